$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that used to sit at the end
#    of the "...ving worse or the same results)" paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the "Treshold 0.1:  NB 46(52), LOG 47(53), J48 45 (54)"
#    paragraph and turn it into
#    "Treshold 0.1 w Wece:  NB 46(52), LOG 47(53), J48 45 (54)"
#    while keeping (and extending) the proof-reading marks.
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Treshold 0\.1:") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $full = $target.Range
    $paraStart = $full.Start

    # Remove the original "Treshold" run (keeps the surrounding
    # spellStart/spellEnd proof marks in place).
    $d.Range($paraStart, $paraStart + 8).Delete()

    # Remove the remainder of the original text (up to, but not
    # including, the paragraph mark).
    $rest = $d.Paragraphs.Item(1)
    $again = $target.Range
    $d.Range($again.Start, $again.End - 1).Delete()

    # Insert the new run/proofErr structure at the start of the
    # (now empty) paragraph.
    $insertionPoint = $d.Range($paraStart, $paraStart)
    $newParaXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:proofErr w:type='spellStart'/>
<w:r><w:t>Treshold</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:t xml:space='preserve'> 0.1</w:t></w:r>
<w:r><w:t xml:space='preserve'> w </w:t></w:r>
<w:proofErr w:type='spellStart'/>
<w:proofErr w:type='gramStart'/>
<w:r><w:t>Wece</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:t>:  NB</w:t></w:r>
<w:proofErr w:type='gramEnd'/>
<w:r><w:t xml:space='preserve'> 46(52), LOG 47(53), J48 45 (54)</w:t></w:r>
</w:p>
"@
    $null = $insertionPoint.InsertXML($newParaXml)
}

# ------------------------------------------------------------------
# 3) Fill the first of the empty paragraphs that follows with the
#    new "Sci : LOG 43,4% (48/3%)" text, re-creating the _GoBack
#    bookmark between "LOG" and " 43,4%".
# ------------------------------------------------------------------
$emptyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Treshold 0\.1 w Wece:") {
        $idx = $p.Index
        $emptyPara = $d.Paragraphs.Item($idx + 1)
        break
    }
}

if ($emptyPara -ne $null) {
    $insertionPoint2 = $d.Range($emptyPara.Range.Start, $emptyPara.Range.Start)
    $sciParaXml = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
<w:proofErr w:type='spellStart'/>
<w:r><w:t>Sci</w:t></w:r>
<w:proofErr w:type='spellEnd'/>
<w:r><w:t xml:space='preserve'> : LOG</w:t></w:r>
<w:bookmarkStart w:id='0' w:name='_GoBack'/>
<w:bookmarkEnd w:id='0'/>
<w:r><w:t xml:space='preserve'> 43,4% (48/3%)</w:t></w:r>
</w:p>
"@
    $null = $insertionPoint2.InsertXML($sciParaXml)
}
